$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.718.33"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "1.621.17"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.33"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5068"
$ws.Range("E6").Value = "  -1.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2556"
$ws.Range("E8").Value = "  -1.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06375"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.23"
$ws.Range("E10").Value = "  -3.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07768"
$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.230"
$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("D13").Value = "1.625.60"
$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").Value = "1.846.06"
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5537"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.46"
$ws.Range("E16").Value = "  -1.89%  "

$ws.Range("D17").Value = "0.0₅7516"
$ws.Range("E17").Value = "  -3.20%  "

$ws.Range("D18").Value = "25.744.77"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.10"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.373"
$ws.Range("E21").Value = "  -1.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.736"
$ws.Range("E22").Value = "  -2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.943"
$ws.Range("E23").Value = "  -2.63%  "

$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.863"
$ws.Range("E25").Value = "  -1.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.44"
$ws.Range("E26").Value = "  -1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1248"
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.711"
$ws.Range("E28").Value = "  -2.53%  "

$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.232"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04852"
$ws.Range("E31").Value = "  -1.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.295"
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.165"
$ws.Range("E33").Value = "  -1.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.538"
$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.361"
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8900"
$ws.Range("E36").Value = "  -3.34%  "

$ws.Range("D37").Value = "1.124.34"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.531"
$ws.Range("E38").Value = "  -2.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5469"
$ws.Range("E39").Value = "  -1.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01553"
$ws.Range("E40").Value = "  -1.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.564"
$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7928"
$ws.Range("E43").Value = "  -2.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.91"
$ws.Range("E44").Value = "  -2.88%  "

$ws.Range("D45").Value = "1.770.39"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("E46").Value = "  -8.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4411"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.51"
$ws.Range("E48").Value = "  -1.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05120"
$ws.Range("E49").Value = "  -3.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.569"
$ws.Range("E50").Value = "  +2.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9966"
$ws.Range("E51").Value = "  -0.89%  "

